# Applies the weekly fruit/vegetable price refresh described in the commit
# "Fruta / hortaliza, semanal" to the "Higo" (fig) price sheet.
#
# The data rows (2-5, 7-18; row 6 is untouched) get new Fecha (date-serial),
# Volumen, Precio minimo/maximo/promedio ponderado and Precio $/Kg values,
# and in a few rows (11-14) the Origen text also changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44302
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 15000
$ws.Range("O2").Value = 15000
$ws.Range("P2").Value = 15000
$ws.Range("S2").Value = 2143

# Row 3
$ws.Range("D3").Value = 44302
$ws.Range("N3").Value = 12000
$ws.Range("O3").Value = 12000
$ws.Range("P3").Value = 12000
$ws.Range("S3").Value = 1714

# Row 4
$ws.Range("D4").Value = 44301
$ws.Range("M4").Value = 100
$ws.Range("N4").Value = 14000
$ws.Range("O4").Value = 14000
$ws.Range("P4").Value = 14000
$ws.Range("S4").Value = 2000

# Row 5
$ws.Range("D5").Value = 44301
$ws.Range("M5").Value = 80
$ws.Range("N5").Value = 12000
$ws.Range("O5").Value = 12000
$ws.Range("P5").Value = 12000
$ws.Range("S5").Value = 1714

# Row 6 is unchanged.

# Row 7
$ws.Range("D7").Value = 44322
$ws.Range("M7").Value = 45
$ws.Range("N7").Value = 12000
$ws.Range("O7").Value = 12000
$ws.Range("P7").Value = 12000
$ws.Range("S7").Value = 1714

# Row 8
$ws.Range("D8").Value = 44322
$ws.Range("N8").Value = 8000
$ws.Range("O8").Value = 8000
$ws.Range("P8").Value = 8000
$ws.Range("S8").Value = 1143

# Row 9
$ws.Range("D9").Value = 44980
$ws.Range("M9").Value = 80

# Row 10
$ws.Range("D10").Value = 44980
$ws.Range("M10").Value = 50
$ws.Range("N10").Value = 13000
$ws.Range("O10").Value = 13000
$ws.Range("P10").Value = 13000
$ws.Range("S10").Value = 1857

# Row 11
$ws.Range("D11").Value = 44320
$ws.Range("M11").Value = 20
$ws.Range("N11").Value = 12000
$ws.Range("O11").Value = 12000
$ws.Range("P11").Value = 12000
$ws.Range("R11").Value = "Región Metropolitana"
$ws.Range("S11").Value = 1714

# Row 12
$ws.Range("D12").Value = 44320
$ws.Range("M12").Value = 30
$ws.Range("N12").Value = 8000
$ws.Range("O12").Value = 8000
$ws.Range("P12").Value = 8000
$ws.Range("R12").Value = "Región Metropolitana"
$ws.Range("S12").Value = 1143

# Row 13
$ws.Range("D13").Value = 44299
$ws.Range("M13").Value = 80
$ws.Range("N13").Value = 15000
$ws.Range("O13").Value = 15000
$ws.Range("P13").Value = 15000
$ws.Range("R13").Value = "Provincia de Santiago"
$ws.Range("S13").Value = 2143

# Row 14
$ws.Range("D14").Value = 44299
$ws.Range("M14").Value = 75
$ws.Range("N14").Value = 12000
$ws.Range("O14").Value = 12000
$ws.Range("P14").Value = 12000
$ws.Range("R14").Value = "Provincia de Santiago"
$ws.Range("S14").Value = 1714

# Row 15
$ws.Range("D15").Value = 44300
$ws.Range("N15").Value = 15000
$ws.Range("O15").Value = 15000
$ws.Range("P15").Value = 15000
$ws.Range("S15").Value = 2143

# Row 16
$ws.Range("D16").Value = 44300

# Row 17
$ws.Range("D17").Value = 44292
$ws.Range("M17").Value = 25
$ws.Range("N17").Value = 16000
$ws.Range("O17").Value = 16000
$ws.Range("P17").Value = 16000
$ws.Range("S17").Value = 2286

# Row 18
$ws.Range("D18").Value = 44292
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 15000
$ws.Range("P18").Value = 15000
$ws.Range("S18").Value = 2143
